$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.195090532302856
$ws.Range("B1").Value = 2.679686069488525
$ws.Range("C1").Value = 9.354107856750488
$ws.Range("D1").Value = 2.075063467025757
$ws.Range("E1").Value = 1.209280729293823
